$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: "Clear" test case becomes "CheckBox" test case -------------
# F20 keeps referencing the same underlying text slot, but the wording
# changes from "2. Clear" to "2. CheckBox".
$ws.Range("F20").Value = "1. Launch Home page`n2. CheckBox"

# G20/H20 description updated to describe the checkbox feature.
$checkBoxDescription = "The user can click a checkbox in the top right hand corner of the map. This will allow an admin to see where the most frequent spots for accidents are, over the past 2 weeks"
$ws.Range("G20").Value = $checkBoxDescription
$ws.Range("H20").Value = $checkBoxDescription

# J20 result flips from Fail to Pass.
$ws.Range("J20").Value = "Pass"

# --- Row 21: new "SearchBar" test case -----------------------------------
# Copy formatting (styles) from row 20 down to row 21 first, then set the
# row height to match the other populated rows, then fill in the values.
$ws.Range("A20:L20").Copy()
$ws.Range("A21:L21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows.Item(21).RowHeight = 85.8

$ws.Range("A21").Value = 1
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = "TC.001"
$ws.Range("D21").Value = "Home"
$ws.Range("E21").Value = "GUI"
$ws.Range("F21").Value = "1. Launch Home page`n2. SearchBar"

$searchBarDescription = "The user can click into the searchbar at the top of the application, where they can enter a location to see if there's any accidents or potholes heading towards their desired destination"
$ws.Range("G21").Value = $searchBarDescription
$ws.Range("H21").Value = $searchBarDescription

$ws.Range("I21").Value = "The user clicks the searchbar, where they enter their desired location to see if there's any accident or potholes along the way"
$ws.Range("J21").Value = "Pass"
$ws.Range("K21").Value = ""
$ws.Range("L21").Value = ""

# I20 actual result now reflects the checkbox behaviour (new text). Set this
# last so the newly-introduced shared string lands after row 21's strings,
# matching the order in which the change was authored.
$ws.Range("I20").Value = "The user clicks the checkbox, and the map will display where all accident markers have been, over the past 2 weeks"

# --- Selection / view -----------------------------------------------------
$ws.Range("F20").Select()
